$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 is intentionally left blank (no data) - matches the source diff's
# empty <row r="44"/> placeholder, so nothing is written there.

# Row 45
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "2025-04-28 02:08:24"
$ws.Cells.Item(45, 3).Value = "James Davis moved Suzuki Battery from shelf space 1 to floor space 2.`nNow James Davis is Frustrated.`n"
$ws.Cells.Item(45, 3).WrapText = $true

# Row 46
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "2025-04-28 02:12:38"
$ws.Cells.Item(46, 3).Value = "James Davis moved New Battery from Ford from floor space 2 to floor space 1.`nNow James Davis is Frustrated.`n"
$ws.Cells.Item(46, 3).WrapText = $true

# The engine auto-expands row height whenever multi-line text lands in a
# cell; AutoFit brings the rows back to the sheet's implicit/default height
# so no stray ht/customHeight attributes are written (matching the diff).
$ws.Rows.Item(45).AutoFit()
$ws.Rows.Item(46).AutoFit()
